$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Mariana Conde"
$ws.Range("B3").Value = 3125048463
$ws.Range("C3").Value = "dsfasdfsadf"

$ws.Range("C3").Select()
